$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.902.55'
$ws.Range('E2').Value = '  -1.77%  '
$ws.Range('D3').Value = '1.812.84'
$ws.Range('E3').Value = '  -0.67%  '
$ws.Range('D4').Value = '0.9993'
$ws.Range('E4').Value = '  -0.21%  '
$ws.Range('D5').Value = '309.88'
$ws.Range('E5').Value = '  -1.29%  '
$ws.Range('D6').Value = '0.9986'
$ws.Range('E6').Value = '  -0.21%  '
$ws.Range('D7').Value = '0.4611'
$ws.Range('E7').Value = '  +3.19%  '
$ws.Range('D8').Value = '0.3743'
$ws.Range('E8').Value = '  -0.32%  '
$ws.Range('D9').Value = '0.07351'
$ws.Range('E9').Value = '  -2.18%  '
$ws.Range('D10').Value = '0.8695'
$ws.Range('D11').Value = '20.45'
$ws.Range('E11').Value = '  -2.95%  '
$ws.Range('D12').Value = '1.738.04'
$ws.Range('E12').Value = '  -4.77%  '
$ws.Range('D13').Value = '5.342'
$ws.Range('E13').Value = '  -1.41%  '
$ws.Range('D14').Value = '6.534'
$ws.Range('E14').Value = '  -3.42%  '
$ws.Range('E15').Value = '  -1.02%  '
$ws.Range('D16').Value = '91.43'
$ws.Range('E16').Value = '  -2.85%  '
$ws.Range('D17').Value = '1.000'
$ws.Range('E17').Value = '  -0.15%  '
$ws.Range('D18').Value = '0.000008722'
$ws.Range('E18').Value = '  -1.06%  '
$ws.Range('D19').Value = '0.9987'
$ws.Range('E19').Value = '  -0.22%  '
$ws.Range('E20').Value = '  -3.01%  '
$ws.Range('D21').Value = '26.911.73'
$ws.Range('E21').Value = '  -1.73%  '
$ws.Range('D22').Value = '5.301'
$ws.Range('D23').Value = '10.74'
$ws.Range('E23').Value = '  -1.82%  '
$ws.Range('D24').Value = '1.984.01'
$ws.Range('E24').Value = '  -3.56%  '
$ws.Range('D25').Value = '1.913'
$ws.Range('E25').Value = '  -3.02%  '
$ws.Range('D26').Value = '151.08'
$ws.Range('E26').Value = '  -0.26%  '
$ws.Range('D27').Value = '18.42'
$ws.Range('E27').Value = '  -0.88%  '
$ws.Range('D28').Value = '2.163'
$ws.Range('E28').Value = '  -9.21%  '
$ws.Range('D29').Value = '5.276'
$ws.Range('E29').Value = '  -1.61%  '
$ws.Range('D30').Value = '114.94'
$ws.Range('E30').Value = '  -2.44%  '
$ws.Range('D31').Value = '0.08886'
$ws.Range('E31').Value = '  +0.61%  '
$ws.Range('D32').Value = '0.7678'
$ws.Range('E32').Value = '  -2.40%  '
$ws.Range('D33').Value = '1.171'
$ws.Range('E33').Value = '  -2.63%  '
$ws.Range('D34').Value = '4.471'
$ws.Range('E34').Value = '  -1.15%  '
$ws.Range('D35').Value = '2.889'
$ws.Range('E35').Value = '  -0.60%  '
$ws.Range('D36').Value = '0.9987'
$ws.Range('E36').Value = '  -0.22%  '
$ws.Range('E37').Value = '  +0.60%  '
$ws.Range('D38').Value = '2.488'
$ws.Range('E38').Value = '  +8.90%  '
$ws.Range('E39').Value = '  -2.07%  '
$ws.Range('D40').Value = '0.05229'
$ws.Range('E40').Value = '  -2.03%  '
$ws.Range('D41').Value = '2.904'
$ws.Range('E41').Value = '  +1.48%  '
$ws.Range('D42').Value = '7.193'
$ws.Range('E42').Value = '  -2.76%  '
$ws.Range('D43').Value = '0.5277'
$ws.Range('E43').Value = '  -0.73%  '
$ws.Range('D44').Value = '0.1661'
$ws.Range('E44').Value = '  -4.06%  '
$ws.Range('D45').Value = '8.582'
$ws.Range('E45').Value = '  -2.04%  '
$ws.Range('D46').Value = '0.5049'
$ws.Range('E46').Value = '  -1.34%  '
$ws.Range('D47').Value = '10.26'
$ws.Range('E47').Value = '  -3.67%  '
$ws.Range('D48').Value = '104.33'
$ws.Range('E48').Value = '  -1.64%  '
$ws.Range('B49').Value = 'NEARProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D49').Value = '1.673'
$ws.Range('E49').Value = '  -1.93%  '
$ws.Range('B50').Value = 'PaxDollar'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D50').Value = '0.9978'
$ws.Range('E50').Value = '  -0.32%  '
$ws.Range('D51').Value = '0.06320'
$ws.Range('E51').Value = '  -0.88%  '
